# Update processing_statistics sheet:
#  - Replace the old Linux path prefix with a new Windows path prefix in column C
#  - Update the "number of rows" counts in column D for the rows whose source files changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "/home/pedro/Documentos/Doutorado/MO412/final_project/Covid19VaccinationProject"
$newPrefix = "C:/Users/Pedro/Documents/Covid19VaccinationProject"

# Row -> new value for column D (number of rows), only for rows whose count changed
$newRowCounts = @{
    2 = 352928
    3 = 1072791
    4 = 234207
    5 = 1019536
    6 = 5084237
    7 = 2645740
    8 = 102000
    9 = 1248800
    10 = 1966481
    11 = 2548231
    12 = 1296558
    13 = 976936
    14 = 7071457
    15 = 2308189
    16 = 1581090
    17 = 3344884
    18 = 3071745
    19 = 1308471
    20 = 2110560
    21 = 1522529
    22 = 3710549
    23 = 536935
    24 = 197379
    25 = 2512839
    26 = 9202153
    27 = 864669
    28 = 746663
    29 = 352928
    30 = 1072791
    31 = 234207
    32 = 1019536
    33 = 5084237
    34 = 2645740
    35 = 102000
    36 = 1248800
    37 = 1966481
    38 = 2548231
    39 = 1296558
    40 = 976936
    41 = 7071457
    42 = 2308189
    43 = 1581090
    44 = 3344884
    45 = 3071745
    46 = 1308471
    47 = 2110560
    48 = 1522529
    49 = 3710549
    50 = 536935
    51 = 197379
    52 = 2512839
    53 = 9202153
    54 = 864669
    55 = 746663
    56 = 1314380
    58 = 1701488
    59 = 1314380
    61 = 1701488
}

$lastRow = 61

for ($row = 2; $row -le $lastRow; $row++) {
    $cCell = $ws.Cells.Item($row, 3)
    $path = $cCell.Value2
    if ($path -ne $null -and $path.ToString().Contains($oldPrefix)) {
        $cCell.Value2 = $path.ToString().Replace($oldPrefix, $newPrefix)
    }

    if ($newRowCounts.ContainsKey($row)) {
        $ws.Cells.Item($row, 4).Value2 = $newRowCounts[$row]
    }
}
